# Allen_Dark_Analysis.xlsx update
# Adds two new test rows to Sheet1:
#   row 17: test13.mht / JEE MAINS 11 April Shift 2 | 2nd Dec,2025
#   row 18: test14.mht / JEE MAINS 04 Apr Shift 2 | 10th Dec,2024
# and extends the colour-scale conditional formatting ranges so they keep
# covering the data through the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new text values column-by-column (matching how the workbook's
# shared string table was originally built) so new unique strings are
# appended to the shared string table in the same relative order as the
# rest of the same-column values.
$ws.Range("B17").Value = "test13.mht"
$ws.Range("B18").Value = "test14.mht"

$ws.Range("C17").Value = "JEE MAINS 11 April Shift 2 | 2nd Dec,2025 | Online Mode"
$ws.Range("C18").Value = "JEE MAINS 04 Apr Shift 2 | 10th Dec,2024 | Online Mode"

$ws.Range("F18").Value = "131 to 151"

# ----- Row 17 numeric values -----
$ws.Range("D17").Value = 225
$ws.Range("G17").Value = 84.28571428571429
$ws.Range("H17").Value = 59
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 5
$ws.Range("K17").Value = 76
$ws.Range("L17").Value = 83.33333333333334
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 4
$ws.Range("O17").Value = 1
$ws.Range("P17").Value = 62
$ws.Range("Q17").Value = 73.91304347826086
$ws.Range("R17").Value = 17
$ws.Range("S17").Value = 6
$ws.Range("T17").Value = 2
$ws.Range("U17").Value = 87
$ws.Range("V17").Value = 95.65217391304348
$ws.Range("W17").Value = 22
$ws.Range("X17").Value = 1
$ws.Range("Y17").Value = 2

# ----- Row 18 numeric values -----
$ws.Range("D18").Value = 268
$ws.Range("E18").Value = 99
$ws.Range("G18").Value = 94.44444444444444
$ws.Range("H18").Value = 68
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 100
$ws.Range("L18").Value = 100
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 77
$ws.Range("Q18").Value = 86.95652173913044
$ws.Range("R18").Value = 20
$ws.Range("S18").Value = 3
$ws.Range("T18").Value = 2
$ws.Range("U18").Value = 91
$ws.Range("V18").Value = 95.83333333333334
$ws.Range("W18").Value = 23
$ws.Range("X18").Value = 1
$ws.Range("Y18").Value = 1

# ----- Extend the colour-scale conditional formatting ranges to row 18 -----
$cfColumns = @("D","E","G","H","I","K","L","M","N","P","Q","R","S","U","V","W","X")
foreach ($col in $cfColumns) {
    $oldRange = $ws.Range("$col`2:$col`16")
    $newRange = $ws.Range("$col`2:$col`18")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}
